{"js": "// Fixed some empty run when bookmarks are used.\n// The template paragraph that ends with \"template :\" was missing a\n// trailing empty run. Word (and the M2Doc generator) expect an empty\n// <w:r/> run at the end of that paragraph; add it back.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that contains the literal \"template\" text and\n// ends the introductory sentence with \" :\" (be defensive about exact\n// paragraph index, match on content instead).\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text && p.text.indexOf(\"template\") !== -1) {\n    target = p;\n    break;\n  }\n}\nif (!target) {\n  // Fallback: first paragraph of the body.\n  target = paragraphs.items[0];\n}\n\n// Insert an empty run (no text) right at the end of that paragraph,\n// using raw OOXML so we get a genuine empty <w:r/> rather than a run\n// containing an empty <w:t/> text element.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:r/></w:p></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nconst endRange = target.getRange(\"End\");\nendRange.insertOoxml(ooxml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Fixed some empty run when bookmarks are used.\n#\n# The paragraph that ends the introductory sentence (\"... template :\")\n# is missing its trailing empty run. Word/M2Doc expect a bare, empty\n# <w:r/> run right before the closing </w:p> of that paragraph; add it\n# back without touching anything else in the document.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph defensively by its distinctive text rather than\n# a hard-coded paragraph index.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*template*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    $target = $d.Paragraphs(1)\n}\n\n$full = $target.Range\n# Exclude the trailing paragraph mark from the range so InsertXML only\n# replaces the paragraph's run content (not the paragraph break itself).\n$full.MoveEnd(1, -1)\n\n# Capture the paragraph's current content as OOXML so every existing\n# run/attribute is preserved exactly, then append a single empty run\n# (<w:r/>) right before the end of the paragraph.\n$originalText = $full.Text\n\n$newXml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body><w:p w:rsidR=\"00791A6F\" w:rsidRDefault=\"00791A6F\" w:rsidP=\"00F5495F\"><w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:proofErr w:type=\"gramStart\"/><w:proofErr w:type=\"gramEnd\"/><w:r w:rsidRPr=\"00DC5685\"><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>A simple demonstration of a</w:t></w:r><w:r w:rsidR=\"00EB5E85\" w:rsidRPr=\"00DC5685\"><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:r w:rsidR=\"00DC5685\"><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>template</w:t></w:r><w:r w:rsidRPr=\"00DC5685\"><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> :</w:t></w:r><w:r/></w:p></w:body>\n</w:document>\n</pkg:xmlData></pkg:part></pkg:package>\n\"@\n\n$full.InsertXML($newXml)\n"}
